$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells per diff ---
# Q58: 1 -> 0
$ws.Cells.Item(58, 17).Value = 0
# Q66: 2 -> 0
$ws.Cells.Item(66, 17).Value = 0
# O929: 0 -> 2
$ws.Cells.Item(929, 15).Value = 2
# R931: inlineStr(empty) -> numeric 0
$ws.Cells.Item(931, 18).Value = 0
# R932: inlineStr(empty) -> numeric 0
$ws.Cells.Item(932, 18).Value = 0

# --- Append new rows 933-951 ---
# Row 933
$ws.Cells.Item(933, 1).Value = 45474
$ws.Cells.Item(933, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(933, 2).Value = 1430.349975585938
$ws.Cells.Item(933, 3).Value = 1498
$ws.Cells.Item(933, 4).Value = 1424.150024414062
$ws.Cells.Item(933, 5).Value = 1461.349975585938
$ws.Cells.Item(933, 6).Value = 1421.995971679688
$ws.Cells.Item(933, 7).Value = 9484383
$ws.Cells.Item(933, 8).Value = 2024
$ws.Cells.Item(933, 9).Value = 7
$ws.Cells.Item(933, 10).Value = 1
$ws.Cells.Item(933, 11).Value = 0
$ws.Cells.Item(933, 12).Value = 0
$ws.Cells.Item(933, 13).Value = 0
$ws.Cells.Item(933, 14).Value = 27
$ws.Cells.Item(933, 15).Value = 0
$ws.Cells.Item(933, 16).Value = 0
$ws.Cells.Item(933, 17).Value = 0

# Row 934
$ws.Cells.Item(934, 1).Value = 45481
$ws.Cells.Item(934, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(934, 2).Value = 1461.349975585938
$ws.Cells.Item(934, 3).Value = 1511.800048828125
$ws.Cells.Item(934, 4).Value = 1441.25
$ws.Cells.Item(934, 5).Value = 1505.050048828125
$ws.Cells.Item(934, 6).Value = 1464.519287109375
$ws.Cells.Item(934, 7).Value = 8440701
$ws.Cells.Item(934, 8).Value = 2024
$ws.Cells.Item(934, 9).Value = 7
$ws.Cells.Item(934, 10).Value = 8
$ws.Cells.Item(934, 11).Value = 0
$ws.Cells.Item(934, 12).Value = 0
$ws.Cells.Item(934, 13).Value = 0
$ws.Cells.Item(934, 14).Value = 28
$ws.Cells.Item(934, 15).Value = 0
$ws.Cells.Item(934, 16).Value = 0
$ws.Cells.Item(934, 17).Value = 0

# Row 935
$ws.Cells.Item(935, 1).Value = 45488
$ws.Cells.Item(935, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(935, 2).Value = 1520.099975585938
$ws.Cells.Item(935, 3).Value = 1545.699951171875
$ws.Cells.Item(935, 4).Value = 1486.5
$ws.Cells.Item(935, 5).Value = 1491.400024414062
$ws.Cells.Item(935, 6).Value = 1451.23681640625
$ws.Cells.Item(935, 7).Value = 9426194
$ws.Cells.Item(935, 8).Value = 2024
$ws.Cells.Item(935, 9).Value = 7
$ws.Cells.Item(935, 10).Value = 15
$ws.Cells.Item(935, 11).Value = 0
$ws.Cells.Item(935, 12).Value = 0
$ws.Cells.Item(935, 13).Value = 0
$ws.Cells.Item(935, 14).Value = 29
$ws.Cells.Item(935, 15).Value = 0
$ws.Cells.Item(935, 16).Value = 0
$ws.Cells.Item(935, 17).Value = 1

# Row 936
$ws.Cells.Item(936, 1).Value = 45495
$ws.Cells.Item(936, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(936, 2).Value = 1479.949951171875
$ws.Cells.Item(936, 3).Value = 1547.550048828125
$ws.Cells.Item(936, 4).Value = 1444.25
$ws.Cells.Item(936, 5).Value = 1541.150024414062
$ws.Cells.Item(936, 6).Value = 1527.43017578125
$ws.Cells.Item(936, 7).Value = 16845547
$ws.Cells.Item(936, 8).Value = 2024
$ws.Cells.Item(936, 9).Value = 7
$ws.Cells.Item(936, 10).Value = 22
$ws.Cells.Item(936, 11).Value = 0
$ws.Cells.Item(936, 12).Value = 0
$ws.Cells.Item(936, 13).Value = 0
$ws.Cells.Item(936, 14).Value = 30
$ws.Cells.Item(936, 15).Value = 0
$ws.Cells.Item(936, 16).Value = 0
$ws.Cells.Item(936, 17).Value = 0

# Row 937
$ws.Cells.Item(937, 1).Value = 45502
$ws.Cells.Item(937, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(937, 2).Value = 1545.650024414062
$ws.Cells.Item(937, 3).Value = 1569
$ws.Cells.Item(937, 4).Value = 1502.300048828125
$ws.Cells.Item(937, 5).Value = 1507.699951171875
$ws.Cells.Item(937, 6).Value = 1494.277954101562
$ws.Cells.Item(937, 7).Value = 9988347
$ws.Cells.Item(937, 8).Value = 2024
$ws.Cells.Item(937, 9).Value = 7
$ws.Cells.Item(937, 10).Value = 29
$ws.Cells.Item(937, 11).Value = 0
$ws.Cells.Item(937, 12).Value = 0
$ws.Cells.Item(937, 13).Value = 0
$ws.Cells.Item(937, 14).Value = 31
$ws.Cells.Item(937, 15).Value = 0
$ws.Cells.Item(937, 16).Value = 0
$ws.Cells.Item(937, 17).Value = 0

# Row 938
$ws.Cells.Item(938, 1).Value = 45509
$ws.Cells.Item(938, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(938, 2).Value = 1470
$ws.Cells.Item(938, 3).Value = 1515.25
$ws.Cells.Item(938, 4).Value = 1426.75
$ws.Cells.Item(938, 5).Value = 1506.699951171875
$ws.Cells.Item(938, 6).Value = 1493.286865234375
$ws.Cells.Item(938, 7).Value = 10709753
$ws.Cells.Item(938, 8).Value = 2024
$ws.Cells.Item(938, 9).Value = 8
$ws.Cells.Item(938, 10).Value = 5
$ws.Cells.Item(938, 11).Value = 0
$ws.Cells.Item(938, 12).Value = 0
$ws.Cells.Item(938, 13).Value = 0
$ws.Cells.Item(938, 14).Value = 32
$ws.Cells.Item(938, 15).Value = 0
$ws.Cells.Item(938, 16).Value = 0
$ws.Cells.Item(938, 17).Value = 0

# Row 939
$ws.Cells.Item(939, 1).Value = 45516
$ws.Cells.Item(939, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(939, 2).Value = 1506.699951171875
$ws.Cells.Item(939, 3).Value = 1589
$ws.Cells.Item(939, 4).Value = 1489.400024414062
$ws.Cells.Item(939, 5).Value = 1585.300048828125
$ws.Cells.Item(939, 6).Value = 1571.187255859375
$ws.Cells.Item(939, 7).Value = 8449032
$ws.Cells.Item(939, 8).Value = 2024
$ws.Cells.Item(939, 9).Value = 8
$ws.Cells.Item(939, 10).Value = 12
$ws.Cells.Item(939, 11).Value = 0
$ws.Cells.Item(939, 12).Value = 0
$ws.Cells.Item(939, 13).Value = 0
$ws.Cells.Item(939, 14).Value = 33
$ws.Cells.Item(939, 15).Value = 0
$ws.Cells.Item(939, 16).Value = 0
$ws.Cells.Item(939, 17).Value = 0

# Row 940
$ws.Cells.Item(940, 1).Value = 45523
$ws.Cells.Item(940, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(940, 2).Value = 1587.849975585938
$ws.Cells.Item(940, 3).Value = 1633.650024414062
$ws.Cells.Item(940, 4).Value = 1575.849975585938
$ws.Cells.Item(940, 5).Value = 1598.400024414062
$ws.Cells.Item(940, 6).Value = 1584.170532226562
$ws.Cells.Item(940, 7).Value = 8602184
$ws.Cells.Item(940, 8).Value = 2024
$ws.Cells.Item(940, 9).Value = 8
$ws.Cells.Item(940, 10).Value = 19
$ws.Cells.Item(940, 11).Value = 0
$ws.Cells.Item(940, 12).Value = 0
$ws.Cells.Item(940, 13).Value = 0
$ws.Cells.Item(940, 14).Value = 34
$ws.Cells.Item(940, 15).Value = 0
$ws.Cells.Item(940, 16).Value = 0
$ws.Cells.Item(940, 17).Value = 0

# Row 941
$ws.Cells.Item(941, 1).Value = 45530
$ws.Cells.Item(941, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(941, 2).Value = 1610
$ws.Cells.Item(941, 3).Value = 1665
$ws.Cells.Item(941, 4).Value = 1609
$ws.Cells.Item(941, 5).Value = 1636.5
$ws.Cells.Item(941, 6).Value = 1621.931396484375
$ws.Cells.Item(941, 7).Value = 9748421
$ws.Cells.Item(941, 8).Value = 2024
$ws.Cells.Item(941, 9).Value = 8
$ws.Cells.Item(941, 10).Value = 26
$ws.Cells.Item(941, 11).Value = 0
$ws.Cells.Item(941, 12).Value = 0
$ws.Cells.Item(941, 13).Value = 0
$ws.Cells.Item(941, 14).Value = 35
$ws.Cells.Item(941, 15).Value = 0
$ws.Cells.Item(941, 16).Value = 0
$ws.Cells.Item(941, 17).Value = 0

# Row 942
$ws.Cells.Item(942, 1).Value = 45537
$ws.Cells.Item(942, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(942, 2).Value = 1636.5
$ws.Cells.Item(942, 3).Value = 1662
$ws.Cells.Item(942, 4).Value = 1608.900024414062
$ws.Cells.Item(942, 5).Value = 1623.25
$ws.Cells.Item(942, 6).Value = 1608.79931640625
$ws.Cells.Item(942, 7).Value = 7720087
$ws.Cells.Item(942, 8).Value = 2024
$ws.Cells.Item(942, 9).Value = 9
$ws.Cells.Item(942, 10).Value = 2
$ws.Cells.Item(942, 11).Value = 0
$ws.Cells.Item(942, 12).Value = 0
$ws.Cells.Item(942, 13).Value = 0
$ws.Cells.Item(942, 14).Value = 36
$ws.Cells.Item(942, 15).Value = 0
$ws.Cells.Item(942, 16).Value = 0
$ws.Cells.Item(942, 17).Value = 0

# Row 943
$ws.Cells.Item(943, 1).Value = 45544
$ws.Cells.Item(943, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(943, 2).Value = 1623.199951171875
$ws.Cells.Item(943, 3).Value = 1659.900024414062
$ws.Cells.Item(943, 4).Value = 1574.75
$ws.Cells.Item(943, 5).Value = 1656.050048828125
$ws.Cells.Item(943, 6).Value = 1641.307373046875
$ws.Cells.Item(943, 7).Value = 10632757
$ws.Cells.Item(943, 8).Value = 2024
$ws.Cells.Item(943, 9).Value = 9
$ws.Cells.Item(943, 10).Value = 9
$ws.Cells.Item(943, 11).Value = 0
$ws.Cells.Item(943, 12).Value = 0
$ws.Cells.Item(943, 13).Value = 0
$ws.Cells.Item(943, 14).Value = 37
$ws.Cells.Item(943, 15).Value = 0
$ws.Cells.Item(943, 16).Value = 0
$ws.Cells.Item(943, 17).Value = 0

# Row 944
$ws.Cells.Item(944, 1).Value = 45551
$ws.Cells.Item(944, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(944, 2).Value = 1651.099975585938
$ws.Cells.Item(944, 3).Value = 1672
$ws.Cells.Item(944, 4).Value = 1582.400024414062
$ws.Cells.Item(944, 5).Value = 1622.050048828125
$ws.Cells.Item(944, 6).Value = 1607.610107421875
$ws.Cells.Item(944, 7).Value = 12507410
$ws.Cells.Item(944, 8).Value = 2024
$ws.Cells.Item(944, 9).Value = 9
$ws.Cells.Item(944, 10).Value = 16
$ws.Cells.Item(944, 11).Value = 0
$ws.Cells.Item(944, 12).Value = 0
$ws.Cells.Item(944, 13).Value = 0
$ws.Cells.Item(944, 14).Value = 38
$ws.Cells.Item(944, 15).Value = 0
$ws.Cells.Item(944, 16).Value = 0
$ws.Cells.Item(944, 17).Value = 0

# Row 945
$ws.Cells.Item(945, 1).Value = 45558
$ws.Cells.Item(945, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(945, 2).Value = 1624.949951171875
$ws.Cells.Item(945, 3).Value = 1664.800048828125
$ws.Cells.Item(945, 4).Value = 1590.050048828125
$ws.Cells.Item(945, 5).Value = 1609.25
$ws.Cells.Item(945, 6).Value = 1594.923950195312
$ws.Cells.Item(945, 7).Value = 11537305
$ws.Cells.Item(945, 8).Value = 2024
$ws.Cells.Item(945, 9).Value = 9
$ws.Cells.Item(945, 10).Value = 23
$ws.Cells.Item(945, 11).Value = 0
$ws.Cells.Item(945, 12).Value = 0
$ws.Cells.Item(945, 13).Value = 0
$ws.Cells.Item(945, 14).Value = 39
$ws.Cells.Item(945, 15).Value = 0
$ws.Cells.Item(945, 16).Value = 0
$ws.Cells.Item(945, 17).Value = 0

# Row 946
$ws.Cells.Item(946, 1).Value = 45565
$ws.Cells.Item(946, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(946, 2).Value = 1606
$ws.Cells.Item(946, 3).Value = 1648.400024414062
$ws.Cells.Item(946, 4).Value = 1567.849975585938
$ws.Cells.Item(946, 5).Value = 1616.449951171875
$ws.Cells.Item(946, 6).Value = 1602.059814453125
$ws.Cells.Item(946, 7).Value = 13897142
$ws.Cells.Item(946, 8).Value = 2024
$ws.Cells.Item(946, 9).Value = 9
$ws.Cells.Item(946, 10).Value = 30
$ws.Cells.Item(946, 11).Value = 0
$ws.Cells.Item(946, 12).Value = 0
$ws.Cells.Item(946, 13).Value = 0
$ws.Cells.Item(946, 14).Value = 40
$ws.Cells.Item(946, 15).Value = 2
$ws.Cells.Item(946, 16).Value = 0
$ws.Cells.Item(946, 17).Value = 0

# Row 947
$ws.Cells.Item(947, 1).Value = 45572
$ws.Cells.Item(947, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(947, 2).Value = 1624.949951171875
$ws.Cells.Item(947, 3).Value = 1675.949951171875
$ws.Cells.Item(947, 4).Value = 1600.199951171875
$ws.Cells.Item(947, 5).Value = 1647.050048828125
$ws.Cells.Item(947, 6).Value = 1632.387451171875
$ws.Cells.Item(947, 7).Value = 7645087
$ws.Cells.Item(947, 8).Value = 2024
$ws.Cells.Item(947, 9).Value = 10
$ws.Cells.Item(947, 10).Value = 7
$ws.Cells.Item(947, 11).Value = 0
$ws.Cells.Item(947, 12).Value = 0
$ws.Cells.Item(947, 13).Value = 0
$ws.Cells.Item(947, 14).Value = 41
$ws.Cells.Item(947, 15).Value = 0
$ws.Cells.Item(947, 16).Value = 0
$ws.Cells.Item(947, 17).Value = 0

# Row 948
$ws.Cells.Item(948, 1).Value = 45579
$ws.Cells.Item(948, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(948, 2).Value = 1648
$ws.Cells.Item(948, 3).Value = 1709.900024414062
$ws.Cells.Item(948, 4).Value = 1645.599975585938
$ws.Cells.Item(948, 5).Value = 1687.900024414062
$ws.Cells.Item(948, 6).Value = 1672.873779296875
$ws.Cells.Item(948, 7).Value = 12534917
$ws.Cells.Item(948, 8).Value = 2024
$ws.Cells.Item(948, 9).Value = 10
$ws.Cells.Item(948, 10).Value = 14
$ws.Cells.Item(948, 11).Value = 0
$ws.Cells.Item(948, 12).Value = 0
$ws.Cells.Item(948, 13).Value = 0
$ws.Cells.Item(948, 14).Value = 42
$ws.Cells.Item(948, 15).Value = 0
$ws.Cells.Item(948, 16).Value = 0
$ws.Cells.Item(948, 17).Value = 0

# Row 949
$ws.Cells.Item(949, 1).Value = 45586
$ws.Cells.Item(949, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(949, 2).Value = 1745
$ws.Cells.Item(949, 3).Value = 1761.849975585938
$ws.Cells.Item(949, 4).Value = 1685.599975585938
$ws.Cells.Item(949, 5).Value = 1716.449951171875
$ws.Cells.Item(949, 6).Value = 1701.169555664062
$ws.Cells.Item(949, 7).Value = 15895552
$ws.Cells.Item(949, 8).Value = 2024
$ws.Cells.Item(949, 9).Value = 10
$ws.Cells.Item(949, 10).Value = 21
$ws.Cells.Item(949, 11).Value = 0
$ws.Cells.Item(949, 12).Value = 0
$ws.Cells.Item(949, 13).Value = 0
$ws.Cells.Item(949, 14).Value = 43
$ws.Cells.Item(949, 15).Value = 0
$ws.Cells.Item(949, 16).Value = 0
$ws.Cells.Item(949, 17).Value = 0

# Row 950
$ws.Cells.Item(950, 1).Value = 45593
$ws.Cells.Item(950, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(950, 2).Value = 1710
$ws.Cells.Item(950, 3).Value = 1727.949951171875
$ws.Cells.Item(950, 4).Value = 1595.550048828125
$ws.Cells.Item(950, 5).Value = 1603.650024414062
$ws.Cells.Item(950, 6).Value = 1589.373779296875
$ws.Cells.Item(950, 7).Value = 7964167
$ws.Cells.Item(950, 8).Value = 2024
$ws.Cells.Item(950, 9).Value = 10
$ws.Cells.Item(950, 10).Value = 28
$ws.Cells.Item(950, 11).Value = 0
$ws.Cells.Item(950, 12).Value = 0
$ws.Cells.Item(950, 13).Value = 0
$ws.Cells.Item(950, 14).Value = 44
$ws.Cells.Item(950, 15).Value = 0
$ws.Cells.Item(950, 16).Value = 0
$ws.Cells.Item(950, 17).Value = 0

# Row 951
$ws.Cells.Item(951, 1).Value = 45600
$ws.Cells.Item(951, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(951, 2).Value = 1614
$ws.Cells.Item(951, 3).Value = 1715.5
$ws.Cells.Item(951, 4).Value = 1608.050048828125
$ws.Cells.Item(951, 5).Value = 1681.349975585938
$ws.Cells.Item(951, 6).Value = 1681.349975585938
$ws.Cells.Item(951, 7).Value = 10455613
$ws.Cells.Item(951, 8).Value = 2024
$ws.Cells.Item(951, 9).Value = 11
$ws.Cells.Item(951, 10).Value = 4
$ws.Cells.Item(951, 11).Value = 0
$ws.Cells.Item(951, 12).Value = 0
$ws.Cells.Item(951, 13).Value = 0
$ws.Cells.Item(951, 14).Value = 45
$ws.Cells.Item(951, 15).Value = 0
$ws.Cells.Item(951, 16).Value = 0
$ws.Cells.Item(951, 17).Value = 0

